$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row below the existing data row, inheriting row 2's formatting.
$ws.Rows("3:3").Insert()

# Extend the styled/used range from column H into columns I and J on rows 2-3,
# copying the formatting that already exists on H2:H3.
$ws.Range("H2:H3").Copy()
$ws.Range("I2:J3").PasteSpecial(-4122)

# Row 2 (first pull) values: cable size / express-local / diameter / weight.
$ws.Range("E2").Value = "7C#14"
$ws.Range("G2").Value = 0.99
$ws.Range("H2").Value = 1.39

# Row 3 (second pull) values.
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "2C#2"

# Merge the per-bundle columns across the two pull rows (not D/E, which vary per pull).
$ws.Range("A2:A3").Merge()
$ws.Range("B2:B3").Merge()
$ws.Range("C2:C3").Merge()
$ws.Range("F2:F3").Merge()
$ws.Range("G2:G3").Merge()
$ws.Range("H2:H3").Merge()
$ws.Range("I2:I3").Merge()
$ws.Range("J2:J3").Merge()
